# Adds the "ODI Batting Extra" and "ODI Bowling Extra" worksheets to the
# workbook, appended after the existing "ODI Bowling" sheet, and populates
# them with per-match batting/bowling "extra" stats.
#
# Cell-encoding convention used by the data tables below: each cell is a
# string "COL|KIND|VALUE" where KIND is "S" for text (written with a
# leading apostrophe so Excel doesn't auto-coerce numeric-looking strings
# like "3560" or "13.54%" into numbers) and "N" for a genuine numeric value.
# (No data value contains "|", so a plain Split on it is unambiguous.)

function Set-RowCells($ws, $rowNum, $cellDefs) {
    foreach ($def in $cellDefs) {
        $parts = $def.Split('|')
        $col = $parts[0]
        $kind = $parts[1]
        $val = $parts[2]
        $colIndex = [int][char]$col - [int][char]'A' + 1
        $cell = $ws.Cells.Item($rowNum, $colIndex)
        if ($kind -eq 'N') {
            $cell.Value = [double]$val
        } else {
            $cell.Value = "'" + $val
        }
    }
}

function Format-HeaderRow($ws, $lastCol) {
    $headerRange = $ws.Range("A1:" + $lastCol + "1")
    $headerRange.Font.Bold = $true
    $headerRange.HorizontalAlignment = -4108
    $headerRange.VerticalAlignment = -4160
    $headerRange.Borders.LineStyle = 1
    $headerRange.Borders.Weight = 2
}


$wb = $excel.ActiveWorkbook

# --- data tables -----------------------------------------------------

$sheet4Rows = @(
    @("A|S|MATCH_CODE","B|S|BATTING_POSITION","C|S|NUM_4","D|S|NUM_6","E|S|PERCENT_RUNS_OF_TOTAL","F|S|MAN_OF_MATCH"),
    @("A|S|3560","B|N|3","C|S|0","D|S|0","F|S|NO"),
    @("A|S|3561","B|N|5","C|S|0","D|S|0","E|S|13.54%","F|S|NO"),
    @("A|S|3632","F|S|NO"),
    @("A|S|3633","B|N|4","C|S|0","D|S|0","E|S|0.99%","F|S|NO"),
    @("A|S|3649","F|S|NO"),
    @("A|S|3650","B|N|4","C|S|0","D|S|0","E|S|4.30%","F|S|NO"),
    @("A|S|3869","B|N|5","C|S|0","D|S|0","E|S|5.17%","F|S|NO"),
    @("A|S|3871","B|N|6","C|S|2","D|S|0","E|S|19.25%","F|S|NO"),
    @("A|S|3873","F|S|NO"),
    @("A|S|3936","B|N|4","C|S|6","D|S|0","E|S|27.91%","F|S|NO"),
    @("A|S|3938","F|S|NO"),
    @("A|S|3941","F|S|NO"),
    @("A|S|3991","F|S|NO"),
    @("A|S|3993","B|N|10","C|S|0","D|S|0","E|S|0.84%","F|S|NO"),
    @("A|S|3996","B|N|5","C|S|1","D|S|0","E|S|10.81%","F|S|NO"),
    @("A|S|4006","B|N|5","C|S|4","D|S|0","E|S|10.96%","F|S|NO"),
    @("A|S|4009","B|N|5","C|S|3","D|S|0","E|S|10.45%","F|S|NO"),
    @("A|S|4189","B|N|5","C|S|3","D|S|0","E|S|23.79%","F|S|NO"),
    @("A|S|4190","B|N|5","C|S|0","D|S|0","E|S|1.10%","F|S|NO"),
    @("A|S|4192","B|N|4","C|S|5","D|S|0","E|S|26.77%","F|S|NO"),
    @("A|S|4195","B|N|5","C|S|2","D|S|0","E|S|14.86%","F|S|NO"),
    @("A|S|4198","B|N|4","C|S|3","D|S|0","E|S|22.75%","F|S|NO"),
    @("A|S|4200","B|N|4","C|S|7","D|S|0","E|S|37.74%","F|S|NO"),
    @("A|S|4202","B|N|4","C|S|5","D|S|0","E|S|28.86%","F|S|NO"),
    @("A|S|4203","F|S|NO"),
    @("A|S|4257","B|N|4","C|S|1","D|S|0","E|S|5.45%","F|S|NO"),
    @("A|S|4259","F|S|NO"),
    @("A|S|4262","F|S|NO"),
    @("A|S|4290","B|N|4","C|S|6","D|S|0","E|S|21.93%","F|S|NO"),
    @("A|S|4299","B|N|5","C|S|1","D|S|0","E|S|8.70%","F|S|NO"),
    @("A|S|4301","B|N|5","C|S|4","D|S|1","E|S|15.41%","F|S|NO"),
    @("A|S|4306","B|N|4","C|S|3","D|S|0","E|S|8.70%","F|S|NO"),
    @("A|S|4309","B|N|4","C|S|0","D|S|0","E|S|2.63%","F|S|NO"),
    @("A|S|4315","B|N|4","C|S|9","D|S|0","E|S|34.30%","F|S|NO"),
    @("A|S|4323","B|N|4","C|S|1","D|S|0","E|S|6.40%","F|S|NO"),
    @("A|S|4326","F|S|NO"),
    @("A|S|4332","F|S|NO"),
    @("A|S|4335","B|N|3","C|S|0","D|S|0","E|S|5.50%","F|S|NO"),
    @("A|S|4340","B|N|3","C|S|0","D|S|0","F|S|NO"),
    @("A|S|4446","B|N|4","C|S|8","D|S|1","E|S|31.54%","F|S|NO"),
    @("A|S|4448","B|N|4","C|S|0","D|S|0","E|S|3.76%","F|S|NO"),
    @("A|S|4525","F|S|NO"),
    @("A|S|4528","B|N|4","C|S|1","D|S|2","E|S|22.78%","F|S|NO"),
    @("A|S|4530","B|N|4","C|S|1","D|S|0","E|S|11.02%","F|S|NO"),
    @("A|S|4537","B|N|4","C|S|3","D|S|1","E|S|13.02%","F|S|NO"),
    @("A|S|4538","F|S|NO"),
    @("A|S|4539","B|N|4","C|S|0","D|S|0","E|S|1.04%","F|S|NO"),
    @("A|S|4582","B|N|4","C|S|13","D|S|0","E|S|31.88%","F|S|NO"),
    @("A|S|4585","B|N|4","C|S|0","D|S|0","E|S|0.44%","F|S|NO"),
    @("A|S|4588","B|N|4","C|S|5","D|S|0","E|S|27.74%","F|S|NO"),
    @("A|S|4671","F|S|NO"),
    @("A|S|4674","B|N|4","C|S|2","D|S|0","E|S|12.28%","F|S|NO"),
    @("A|S|4675","F|S|NO")
)

$sheet5Rows = @(
    @("A|S|MATCH_CODE","B|S|MAIDEN_OVERS","C|S|PERCENT_WICKETS_OF_ALL"),
    @("A|S|3633","B|S|0"),
    @("A|S|3991")
)

# --- add "ODI Batting Extra" after the last existing sheet -----------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$rowNum = 1
foreach ($row in $sheet4Rows) {
    Set-RowCells $battingExtra $rowNum $row
    $rowNum++
}
Format-HeaderRow $battingExtra "F"

# --- add "ODI Bowling Extra" after "ODI Batting Extra" ----------------

$bowlingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"

$rowNum = 1
foreach ($row in $sheet5Rows) {
    Set-RowCells $bowlingExtra $rowNum $row
    $rowNum++
}
Format-HeaderRow $bowlingExtra "C"
